# Update column F ("dSF") values for specific rows, re-pulling / recalculating
# the mean-adjusted figures as described in the commit message
# ("repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = 2
    16 = -3
    26 = -6
    27 = -4
    28 = 8
    39 = 2
    42 = 1
    43 = 0
    54 = -1
    55 = 2
    58 = 3
    74 = -6
    77 = -2
    78 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
